$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = 18

$ws.Range("A13").Value = "rdfdata.org"
$ws.Range("B13").Value = 6

$ws.Range("A14").Value = "swoogle.umbc.edu(counting)"
$ws.Range("B14").Value = 1202
$ws.Range("C14").Value = 116

$ws.Range("C14").Select()
